$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Badan Kesatuan Bangsa Dan Politik"
$ws.Range("C6").Value = "Pengembangan Smart City Melalui JAKI"
$ws.Range("D6").Value = 18
$ws.Range("E6").Value = 12
$ws.Range("F6").Value = 24
$ws.Range("G6").Value = 15
$ws.Range("H6").Value = 69
$ws.Range("I6").Value = "bagus"
